# "1st changes of mifos to finflux"
#
# Repayment schedule sheet gains a new (blank-header) column between the
# existing "Late" column and the "# / Principal" columns before it, pushing
# Late / heading(date) / Outstanding one column to the right (N->O, O->P,
# P->Q). Active-sheet/selection bookkeeping also moves from the
# "Transactions" tab back onto "Repayment schedule".

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsTransactions = $wb.Worksheets.Item("Transactions")

# Leave the old selection on "Transactions" where the diff shows it ending up,
# before switching the active tab away from it.
$wsTransactions.Range("F19").Select() | Out-Null

# Insert a new blank column before the old column N ("Late"); this shifts
# "Late" (N->O), the blank-heading date column (O->P) and "Outstanding"
# (P->Q) one slot to the right, exactly matching the diff.
$wsSchedule.Columns("N").Insert()

# The inherited column width on the newly inserted column matches the
# neighbouring "Principal" column (M).
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# "Repayment schedule" becomes the active sheet/tab with the new selection.
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("S7").Select() | Out-Null
